# modify offset for nodes labels
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("for_circular_layout")
$wsRegions = $wb.Worksheets.Item("Regions")

# --- Update the C/D offset values on "for_circular_layout" ---
$wsData.Range("D21").Value = 85
$wsData.Range("D22").Value = 120
$wsData.Range("D23").Value = 80
$wsData.Range("C24").Value = -30
$wsData.Range("C25").Value = -40
$wsData.Range("D25").Value = 95
$wsData.Range("C26").Value = -45
$wsData.Range("D26").Value = 95
$wsData.Range("D27").Value = 90
$wsData.Range("D28").Value = 45
$wsData.Range("D29").Value = 70
$wsData.Range("C32").Value = -75
$wsData.Range("D33").Value = 25
$wsData.Range("C34").Value = -70
$wsData.Range("D34").Value = 35
$wsData.Range("C36").Value = -115
$wsData.Range("C37").Value = -115
$wsData.Range("C38").Value = -95
$wsData.Range("C39").Value = -95
$wsData.Range("C40").Value = -90
$wsData.Range("C42").Value = -95
$wsData.Range("C43").Value = -65
$wsData.Range("C44").Value = -75
$wsData.Range("C45").Value = -85
$wsData.Range("C46").Value = -55
$wsData.Range("C47").Value = -95
$wsData.Range("C48").Value = -45
$wsData.Range("D48").Value = -45
$wsData.Range("C49").Value = -65
$wsData.Range("D49").Value = -70
$wsData.Range("C52").Value = -40
$wsData.Range("D52").Value = -75
$wsData.Range("C53").Value = -30
$wsData.Range("D53").Value = -65
$wsData.Range("C54").Value = -40
$wsData.Range("D54").Value = -95
$wsData.Range("C56").Value = -20
$wsData.Range("D56").Value = -90
$wsData.Range("C57").Value = -10
$wsData.Range("D57").Value = -95
$wsData.Range("D58").Value = -15
$wsData.Range("D59").Value = -15
$wsData.Range("D60").Value = -15
$wsData.Range("C61").Value = 0
$wsData.Range("D61").Value = -15
$wsData.Range("C62").Value = 0
$wsData.Range("D62").Value = -15
$wsData.Range("C63").Value = 0
$wsData.Range("D63").Value = -15
$wsData.Range("D64").Value = -15
$wsData.Range("C65").Value = 5
$wsData.Range("D65").Value = -15
$wsData.Range("C66").Value = 5
$wsData.Range("D66").Value = -15
$wsData.Range("C67").Value = 5
$wsData.Range("D67").Value = -15
$wsData.Range("C68").Value = 5
$wsData.Range("C69").Value = 10
$wsData.Range("C70").Value = 10
$wsData.Range("C71").Value = 10
$wsData.Range("C72").Value = 10
$wsData.Range("C73").Value = 10
$wsData.Range("C74").Value = 10
$wsData.Range("C75").Value = 10
$wsData.Range("C76").Value = 10
$wsData.Range("C77").Value = 10

# D59's fill (highlight) is restyled to match its neighbour D58's highlighted fill
$wsData.Range("D58").Copy()
$wsData.Range("D59").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Sheet view / selection bookkeeping ---
# "for_circular_layout" becomes the active tab, with D67 selected
$wsData.Activate()
$wsData.Range("D67").Select()

# "Regions" sheet loses its tabSelected flag, but keeps its scroll position
$wsRegions.Range("A77").Select()

# Restore the data sheet as the frontmost/active sheet and set the
# window size recorded for the workbook view
$wsData.Activate()
$excel.ActiveWindow.Width = 17505
$excel.ActiveWindow.Height = 12270
